$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the row header from "A" to "Server Room"
$ws.Range("A4").Value = "Server Room"

# Shorten the title (drop the "(Area Server Room)" suffix)
$ws.Range("A1").Value = "SUBNETTING RETE DI CLASSE C"

# Move the active selection to I12, matching the saved view state
$ws.Range("I12").Select()
